$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in B (Year), C (IVPR Sl. No.), D (Page) for rows that were
# previously blank in those three columns; column A is left untouched.
$data = @(
    @(1098, 1995, 834, 397),
    @(1099, 1995, 4692, 597),
    @(1100, 1995, 1174, 415),
    @(1101, 1996, 508, 539),
    @(1102, 1996, 738, 551),
    @(1103, 1996, 548, 541),
    @(1104, 1996, 648, 546),
    @(1105, 1996, 1271, 579),
    @(1106, 1996, 503, 539),
    @(1107, 1996, 2248, 629),
    @(1108, 1996, 864, 557),
    @(1111, 1998, 908, 251),
    @(1113, 1998, 2416, 352),
    @(1114, 1998, 1346, 280),
    @(1115, 1998, 545, 227),
    @(1116, 1998, 913, 251),
    @(1117, 1998, 5994, 590),
    @(1119, 1998, 5991, 590),
    @(1120, 1999, 497, 195),
    @(1121, 1999, 324, 185),
    @(1122, 1999, 1211, 237),
    @(1123, 1999, 704, 207),
    @(1124, 1999, 264, 181),
    @(1125, 1999, 1210, 237),
    @(1126, 1999, 287, 182),
    @(1127, 1999, 686, 206),
    @(1128, 1999, 831, 214),
    @(1129, 2000, 1829, 406),
    @(1130, 2000, 667, 300),
    @(1131, 2001, 683, 262),
    @(1132, 2001, 485, 246),
    @(1133, 2001, 2106, 381),
    @(1134, 2001, 399, 239),
    @(1135, 2001, 544, 251),
    @(1136, 2001, 1375, 320),
    @(1137, 2001, 1508, 331),
    @(1138, 2002, 569, 208),
    @(1139, 2002, 320, 184),
    @(1140, 2002, 37, 155),
    @(1141, 2003, 203, 121),
    @(1142, 2003, 376, 131),
    @(1143, 2003, 1123, 176),
    @(1144, 2003, 1047, 172),
    @(1145, 2003, 303, 127),
    @(1146, 2003, 356, 130),
    @(1147, 2003, 239, 123),
    @(1148, 2003, 353, 130),
    @(1149, 2004, 270, 142),
    @(1150, 2004, 1074, 202),
    @(1151, 2004, 289, 143),
    @(1152, 2004, 205, 137),
    @(1153, 2005, 141, 101),
    @(1154, 2006, 446, 345),
    @(1155, 2006, 3628, 634),
    @(1156, 2007, 287, 304),
    @(1157, 2007, 1937, 469),
    @(1158, 2007, 310, 307)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Replicate the view-state side effect captured in the saved file: the
# author scrolled down to the newly-edited rows, left the window zoom at
# 100%, and left the selection on B1159.
$excel.ActiveWindow.Zoom = 100
$excel.ActiveWindow.ScrollRow = 1156
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B1159").Select()
